# Data update from DGS's 2021/10/18 report.
# Append one new time-series row (row 95) to Sheet1: date in column A
# (kept as text, matching the existing column-A string entries) and four
# numeric indicators in columns B-E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 95

# Force the date to be stored as literal text (not auto-converted to a
# date serial) while still ending up tagged with the same "yyyy/mm/dd"
# display format the rest of column A already uses.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2021/10/18"
$ws.Cells.Item($row, 1).NumberFormat = "yyyy/mm/dd"

$ws.Cells.Item($row, 2).Value = 84.3
$ws.Cells.Item($row, 3).Value = 84.7
$ws.Cells.Item($row, 4).Value = 1.01
$ws.Cells.Item($row, 5).Value = 1.02

# Move the active selection down to the next empty row, matching the
# cursor position left behind after typing a new row of data.
$ws.Cells.Item($row + 1, 1).Select() | Out-Null
